$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is a label/value table (A=label, B=value, C=value-highlighted).
# In the "before" workbook the B/C value column was off by one row versus
# its label in A for the "Objetivos" block onward. Fix this by inserting a
# new row at 13 (shifting the old rows 13-23 down to 14-24) and then
# re-pointing every label's B/C value to the correct text - including a
# few brand-new value strings that didn't exist before.

# 1) Insert a blank row at 13; this pushes rows 13..23 down to 14..24 and
#    keeps per-cell styles/row heights intact (A/B/C keep styles 1/2/3).
$ws.Rows.Item(13).Insert()

# 2) The new row 13 should look like rows such as row 1 (only B/C filled,
#    no A label). Copy just the number/cell formatting from row 14 (B/C)
#    into row 13 so the new cells carry the correct styles (2 and 3)
#    without pulling in row 14's values, then drop the stray A13 cell that
#    Insert() stamped with column A's style.
$ws.Range("B14:C14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("A13").Clear()

# 3) Write the correct value text for every label now that rows are
#    properly aligned.
$ws.Range("B10").Value = "Possibilitar aos alunos a realização de trabalho de síntese e integração dos conhecimentos adquiridos ao longo do curso, conforme projeto aprovado na disciplina de Trabalho de Conclusão do Curso I."
$ws.Range("C10").Value = "Possibilitar aos alunos a realização de trabalho de síntese e integração dos conhecimentos adquiridos ao longo do curso, conforme projeto aprovado na disciplina de Trabalho de Conclusão do Curso I."

$ws.Range("B13").Value = "1285870 - Marcos Villela Barcza"
$ws.Range("C13").Value = "1285870 - Marcos Villela Barcza"

$ws.Range("B14").Value = "Desenvolvimento do trabalho de conclusão de curso, sob orientação de um professor orientador, o qual deve constituir-se num projeto de engenharia química."
$ws.Range("C14").Value = "Desenvolvimento do trabalho de conclusão de curso, sob orientação de um professor orientador, o qual deve constituir-se num projeto de engenharia química."

$ws.Range("B16").Value = "Elaboração de uma monografia de conclusão de curso que apresente: (1) o tema e sua importância, (2) os objetivos, (3) a revisão bibliográfica, (4) a metodologia científica (5) o desenvolvimento do projeto, (6) a análise e discussão dos resultados, (7) as conclusões e (8) referências bibliográficas."
$ws.Range("C16").Value = "Elaboração de uma monografia de conclusão de curso que apresente: (1) o tema e sua importância, (2) os objetivos, (3) a revisão bibliográfica, (4) a metodologia científica (5) o desenvolvimento do projeto, (6) a análise e discussão dos resultados, (7) as conclusões e (8) referências bibliográficas."

$ws.Range("B19").Value = "Reuniões periódicas com o orientador e realização do trabalho de conclusão de curso conforme orientação e apresentação de uma monografia final, conforme norma do Departamento de Engenharia Química."
$ws.Range("C19").Value = "Reuniões periódicas com o orientador e realização do trabalho de conclusão de curso conforme orientação e apresentação de uma monografia final, conforme norma do Departamento de Engenharia Química."

$ws.Range("B20").Value = "Avaliação da monografia perante uma banca examinadora composta por 3 (três) membros, obrigatoriamente docentes da Escola de Engenharia de Lorena (EEL)."
$ws.Range("C20").Value = "Avaliação da monografia perante uma banca examinadora composta por 3 (três) membros, obrigatoriamente docentes da Escola de Engenharia de Lorena (EEL)."

$ws.Range("B21").Value = "Reapresentação da monografia, preferencialmente para a mesma banca, com as modificações sugeridas para uma nova avaliação."
$ws.Range("C21").Value = "Reapresentação da monografia, preferencialmente para a mesma banca, com as modificações sugeridas para uma nova avaliação."

$ws.Range("B22").Value = "Recomendada pelo orientador"
$ws.Range("C22").Value = "Recomendada pelo orientador"
